$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.197986577181208
$ws.Range("C2").Value = 0.5469798657718121
$ws.Range("J2").Value = 0.006711409395973154
$ws.Range("P2").Value = 0.1342281879194631
$ws.Range("S2").Value = 0.1140939597315436
$ws.Range("B3").Value = 0.01176470588235294
$ws.Range("C3").Value = 0.01764705882352941
$ws.Range("J3").Value = 0.04705882352941176
$ws.Range("P3").Value = 0.7764705882352941
$ws.Range("S3").Value = 0.1470588235294118
$ws.Range("J4").Value = 0.0303030303030303
$ws.Range("P4").Value = 0.696969696969697
$ws.Range("S4").Value = 0.2727272727272727
$ws.Range("B6").Value = 0.05263157894736842
$ws.Range("D6").Value = 0.01578947368421053
$ws.Range("E6").Value = 0.005263157894736842
$ws.Range("F6").Value = 0.06315789473684211
$ws.Range("J6").Value = 0.2421052631578947
$ws.Range("O6").Value = 0.04210526315789474
$ws.Range("Q6").Value = 0.1736842105263158
$ws.Range("R6").Value = 0.07368421052631578
$ws.Range("S6").Value = 0.331578947368421
$ws.Range("B7").Value = 0.1147540983606557
$ws.Range("D7").Value = 0.01092896174863388
$ws.Range("F7").Value = 0.02185792349726776
$ws.Range("J7").Value = 0.1147540983606557
$ws.Range("O7").Value = 0.0273224043715847
$ws.Range("Q7").Value = 0.1584699453551913
$ws.Range("R7").Value = 0.1038251366120219
$ws.Range("S7").Value = 0.4480874316939891
$ws.Range("B8").Value = 0.0796812749003984
$ws.Range("D8").Value = 0.02191235059760956
$ws.Range("F8").Value = 0.05378486055776893
$ws.Range("J8").Value = 0.1135458167330677
$ws.Range("O8").Value = 0.01394422310756972
$ws.Range("Q8").Value = 0.1693227091633466
$ws.Range("R8").Value = 0.09362549800796813
$ws.Range("S8").Value = 0.4541832669322709
$ws.Range("B9").Value = 0.1484375
$ws.Range("D9").Value = 0.0234375
$ws.Range("F9").Value = 0.046875
$ws.Range("J9").Value = 0.1328125
$ws.Range("O9").Value = 0.0234375
$ws.Range("Q9").Value = 0.1484375
$ws.Range("R9").Value = 0.1328125
$ws.Range("S9").Value = 0.34375
$ws.Range("B10").Value = 0.1271777003484321
$ws.Range("D10").Value = 0.01306620209059233
$ws.Range("F10").Value = 0.0627177700348432
$ws.Range("J10").Value = 0.10801393728223
$ws.Range("O10").Value = 0.01393728222996516
$ws.Range("Q10").Value = 0.2116724738675958
$ws.Range("R10").Value = 0.08710801393728224
$ws.Range("S10").Value = 0.3763066202090593
$ws.Range("G11").Value = 0.1832797427652733
$ws.Range("J11").Value = 0.1157556270096463
$ws.Range("K11").Value = 0.2347266881028939
$ws.Range("L11").Value = 0.4598070739549839
$ws.Range("S11").Value = 0.006430868167202572
$ws.Range("G12").Value = 0.7142857142857143
$ws.Range("J12").Value = 0.2312925170068027
$ws.Range("K12").Value = 0.02040816326530612
$ws.Range("L12").Value = 0.0272108843537415
$ws.Range("S12").Value = 0.006802721088435374
$ws.Range("G13").Value = 0.6428571428571429
$ws.Range("J13").Value = 0.3095238095238095
$ws.Range("S13").Value = 0.04761904761904762
$ws.Range("F15").Value = 0.03191489361702127
$ws.Range("H15").Value = 0.1968085106382979
$ws.Range("I15").Value = 0.05319148936170213
$ws.Range("J15").Value = 0.3031914893617021
$ws.Range("K15").Value = 0.09574468085106383
$ws.Range("O15").Value = 0.0425531914893617
$ws.Range("S15").Value = 0.2765957446808511
$ws.Range("F16").Value = 0.02590673575129534
$ws.Range("H16").Value = 0.2020725388601036
$ws.Range("I16").Value = 0.06735751295336788
$ws.Range("J16").Value = 0.3316062176165803
$ws.Range("K16").Value = 0.1295336787564767
$ws.Range("M16").Value = 0.0310880829015544
$ws.Range("O16").Value = 0.06735751295336788
$ws.Range("S16").Value = 0.1450777202072539
$ws.Range("F17").Value = 0.01728395061728395
$ws.Range("H17").Value = 0.1851851851851852
$ws.Range("I17").Value = 0.0691358024691358
$ws.Range("J17").Value = 0.4098765432098765
$ws.Range("K17").Value = 0.1185185185185185
$ws.Range("M17").Value = 0.01728395061728395
$ws.Range("N17").Value = 0.002469135802469136
$ws.Range("O17").Value = 0.04691358024691358
$ws.Range("S17").Value = 0.1333333333333333
$ws.Range("F18").Value = 0.03061224489795918
$ws.Range("H18").Value = 0.1938775510204082
$ws.Range("I18").Value = 0.06122448979591837
$ws.Range("J18").Value = 0.3979591836734694
$ws.Range("K18").Value = 0.08673469387755102
$ws.Range("M18").Value = 0.02551020408163265
$ws.Range("O18").Value = 0.0663265306122449
$ws.Range("S18").Value = 0.1377551020408163
$ws.Range("F19").Value = 0.01157981803143093
$ws.Range("H19").Value = 0.2622001654259719
$ws.Range("I19").Value = 0.05376344086021505
$ws.Range("J19").Value = 0.3655913978494624
$ws.Range("K19").Value = 0.1033912324234905
$ws.Range("M19").Value = 0.02150537634408602
$ws.Range("O19").Value = 0.06286186931348221
$ws.Range("S19").Value = 0.119106699751861

Write-Output "Applied 107 cell updates"
